# Auto-generated: apply scheduled-runner profit/price updates across Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 68.40000000000001
$ws.Range("I6").Value = 14
$ws.Range("K6").Value = 42
$ws.Range("M6").Value = 70

$ws.Range("H28").Value = 8714.25
$ws.Range("I28").Value = 3571.389
$ws.Range("K28").Value = 3571.389
$ws.Range("M28").Value = -3086.389

$ws.Range("H94").Value = 10250.5
$ws.Range("I94").Value = 10562.077
$ws.Range("J94").Value = 6200
$ws.Range("K94").Value = 10562.077
$ws.Range("L94").Value = 6200
$ws.Range("M94").Value = -10111.077
$ws.Range("N94").Value = -7102

$ws.Range("H103").Value = 3518.375
$ws.Range("I103").Value = 554.5714
$ws.Range("J103").Value = 5823.5557
$ws.Range("K103").Value = 1663.7142
$ws.Range("L103").Value = 17470.6671
$ws.Range("M103").Value = -1077.7142
$ws.Range("N103").Value = -18642.6671

$ws.Range("H107").Value = 209.25
$ws.Range("I107").Value = 126.47059
$ws.Range("J107").Value = 678.3333
$ws.Range("K107").Value = 126.47059
$ws.Range("L107").Value = 678.3333
$ws.Range("M107").Value = 1793.52941
$ws.Range("N107").Value = -4518.3333

$ws.Range("H111").Value = 1581.341
$ws.Range("I111").Value = 1226.6111
$ws.Range("K111").Value = 3679.8333
$ws.Range("M111").Value = -612.8333000000002

$ws.Range("H115").Value = 1609
$ws.Range("I115").Value = 939.9091
$ws.Range("J115").Value = 5289
$ws.Range("K115").Value = 2819.7273
$ws.Range("L115").Value = 15867
$ws.Range("M115").Value = -1252.7273
$ws.Range("N115").Value = -19001

$ws.Range("H116").Value = 4057.45
$ws.Range("I116").Value = 3750.7058
$ws.Range("K116").Value = 3750.7058
$ws.Range("M116").Value = -308.7058000000002

$ws.Range("H138").Value = 3257.158
$ws.Range("J138").Value = 3492.875
$ws.Range("L138").Value = 10478.625
$ws.Range("N138").Value = -20758.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2956.5
$ws.Range("I2").Value = 2624.6667
$ws.Range("J2").Value = 3620.1667
$ws.Range("K2").Value = 2624.6667
$ws.Range("L2").Value = 3620.1667
$ws.Range("M2").Value = -2511.6667
$ws.Range("N2").Value = -3846.1667

$ws.Range("H32").Value = 12612.617
$ws.Range("I32").Value = 8101.0713
$ws.Range("K32").Value = 8101.0713
$ws.Range("M32").Value = -7814.0713

$ws.Range("H61").Value = 4813.4
$ws.Range("I61").Value = 4348.3887
$ws.Range("K61").Value = 4348.3887
$ws.Range("M61").Value = -4136.3887

$ws.Range("H63").Value = 2999.8333
$ws.Range("I63").Value = 3099.8
$ws.Range("K63").Value = 3099.8
$ws.Range("M63").Value = -2413.8

$ws.Range("H66").Value = 2999.8333
$ws.Range("I66").Value = 3099.8
$ws.Range("K66").Value = 15499
$ws.Range("M66").Value = -12067

$ws.Range("H88").Value = 7500
$ws.Range("I88").Value = 7500
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 7500
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -7094
$ws.Range("N88").Value = ""

$ws.Range("H91").Value = 7500
$ws.Range("I91").Value = 7500
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 7500
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -6096
$ws.Range("N91").Value = ""

$ws.Range("H116").Value = 2956.5
$ws.Range("I116").Value = 2624.6667
$ws.Range("J116").Value = 3620.1667
$ws.Range("K116").Value = 2624.6667
$ws.Range("L116").Value = 3620.1667
$ws.Range("M116").Value = -330.6667000000002
$ws.Range("N116").Value = -8208.1667

$ws.Range("H136").Value = 4813.4
$ws.Range("I136").Value = 4348.3887
$ws.Range("K136").Value = 13045.1661
$ws.Range("M136").Value = -10495.1661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2956.5
$ws.Range("I3").Value = 2624.6667
$ws.Range("J3").Value = 3620.1667
$ws.Range("K3").Value = 2624.6667
$ws.Range("L3").Value = 3620.1667
$ws.Range("M3").Value = -2510.6667
$ws.Range("N3").Value = -3848.1667

$ws.Range("H22").Value = 348.15384
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -746

$ws.Range("H134").Value = 10998.5
$ws.Range("I134").Value = 9998
$ws.Range("K134").Value = 29994
$ws.Range("M134").Value = -27459

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 3708
$ws.Range("J15").Value = 5474.5
$ws.Range("L15").Value = 5474.5
$ws.Range("N15").Value = -5814.5

$ws.Range("H25").Value = 6641.1665
$ws.Range("I25").Value = 4962.25
$ws.Range("K25").Value = 4962.25
$ws.Range("M25").Value = -4788.25

$ws.Range("H31").Value = 8545.549999999999
$ws.Range("J31").Value = 15999.5
$ws.Range("L31").Value = 15999.5
$ws.Range("N31").Value = -16589.5

$ws.Range("H34").Value = 8545.549999999999
$ws.Range("J34").Value = 15999.5
$ws.Range("L34").Value = 15999.5
$ws.Range("N34").Value = -16403.5

$ws.Range("H58").Value = 2000
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797

$ws.Range("H100").Value = 99999
$ws.Range("J100").Value = 99999
$ws.Range("L100").Value = 99999
$ws.Range("N100").Value = -102163

$ws.Range("H107").Value = 544.2222
$ws.Range("I107").Value = 506.85715
$ws.Range("K107").Value = 506.85715
$ws.Range("M107").Value = 1413.14285

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1030.1765
$ws.Range("I5").Value = 1032.0625
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 3096.1875
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -2984.1875
$ws.Range("N5").Value = -3224

$ws.Range("H17").Value = 318
$ws.Range("I17").Value = 434.2857
$ws.Range("J17").Value = 46.666668
$ws.Range("K17").Value = 1302.8571
$ws.Range("L17").Value = 140.000004
$ws.Range("M17").Value = -1133.8571
$ws.Range("N17").Value = -478.000004

$ws.Range("H104").Value = 10000
$ws.Range("I104").Value = 10000
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("M104").Value = -27379
$ws.Range("N104").Value = -35242

$ws.Range("H135").Value = 1030.1765
$ws.Range("I135").Value = 1032.0625
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9288.5625
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6753.5625
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""

$ws.Range("H102").Value = 1060.0526
$ws.Range("I102").Value = 595.2353000000001
$ws.Range("J102").Value = 5011
$ws.Range("K102").Value = 595.2353000000001
$ws.Range("L102").Value = 5011
$ws.Range("M102").Value = 1026.7647
$ws.Range("N102").Value = -8255

$ws.Range("H107").Value = 433
$ws.Range("I107").Value = 399.5
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 399.5
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1520.5
$ws.Range("N107").Value = -4340

$ws.Range("H113").Value = 3435.087
$ws.Range("I113").Value = 1000.4667
$ws.Range("K113").Value = 1000.4667
$ws.Range("M113").Value = 1169.5333

$ws.Range("H126").Value = 3222.5
$ws.Range("I126").Value = 3222.5
$ws.Range("K126").Value = 9667.5
$ws.Range("M126").Value = -7197.5

$ws.Range("H132").Value = 67036.52
$ws.Range("I132").Value = 87764.28999999999
$ws.Range("K132").Value = 263292.87
$ws.Range("M132").Value = -260762.87

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4975.4287
$ws.Range("I61").Value = 4248.857
$ws.Range("J61").Value = 6428.5713
$ws.Range("K61").Value = 4248.857
$ws.Range("L61").Value = 6428.5713
$ws.Range("M61").Value = -4046.857
$ws.Range("N61").Value = -6832.5713

$ws.Range("H113").Value = 4975.4287
$ws.Range("I113").Value = 4248.857
$ws.Range("J113").Value = 6428.5713
$ws.Range("K113").Value = 4248.857
$ws.Range("L113").Value = 6428.5713
$ws.Range("M113").Value = -2078.857
$ws.Range("N113").Value = -10768.5713

$ws.Range("H132").Value = 9535.182000000001
$ws.Range("J132").Value = 11997.667
$ws.Range("L132").Value = 35993.001
$ws.Range("N132").Value = -41053.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 507.9524
$ws.Range("I113").Value = 455.14285
$ws.Range("K113").Value = 1365.42855
$ws.Range("M113").Value = 804.5714499999999
